# Updated symbol list on Thu Dec 22 05:30:20 UTC 2022 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values on sheet1.
# D-column values are stored as text in the workbook, so values are
# prefixed with an apostrophe to force Excel to keep them as text
# instead of auto-converting the numeric-looking strings to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.24"
$ws.Range("D3").Value = "'22.64"
$ws.Range("D4").Value = "'5.290"
$ws.Range("D5").Value = "'0.05727"
$ws.Range("D6").Value = "'3.440"
$ws.Range("D7").Value = "'0.8096"
$ws.Range("D8").Value = "'0.8678"
$ws.Range("D9").Value = "'0.1431"
$ws.Range("D10").Value = "'0.07336"
$ws.Range("D11").Value = "'0.03066"
$ws.Range("D12").Value = "'0.03127"
$ws.Range("D13").Value = "'0.09398"
$ws.Range("D14").Value = "'3.903"
$ws.Range("D15").Value = "'0.001576"
$ws.Range("D16").Value = "'0.04816"
$ws.Range("D17").Value = "'0.0005841"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("D18").Value = "'0.006148"
$ws.Range("D19").Value = "'0.005126"
$ws.Range("D20").Value = "'0.0009975"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D22").Value = "'3.733"
$ws.Range("D23").Value = "'6.323"
$ws.Range("D40").Value = "'0.03936"
$ws.Range("D41").Value = "'0.006760"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").Value = "'0.1066"
$ws.Range("D44").Value = "'0.008156"
$ws.Range("D45").Value = "'0.00005623"
$ws.Range("D47").Value = "'0.6001"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Write-Host "Updated symbol list values"
